# Auto-generated Excel COM-interop script
# Applies numeric odds updates to Sheet1 per the target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 8.5
$ws.Range("Z2").Value = 8.5
$ws.Range("G3").Value = 2.3
$ws.Range("I3").Value = 3.1
$ws.Range("AF6").Value = 37
$ws.Range("AH6").Value = 150
$ws.Range("I6").Value = 6.5
$ws.Range("L6").Value = 1.37
$ws.Range("Y6").Value = 35
$ws.Range("AA7").Value = 7.2
$ws.Range("AB7").Value = 16
$ws.Range("AC7").Value = 75
$ws.Range("AD7").Value = 600
$ws.Range("AE7").Value = 13
$ws.Range("AF7").Value = 27
$ws.Range("AG7").Value = 15.5
$ws.Range("AI7").Value = 50
$ws.Range("G7").Value = 1.65
$ws.Range("H7").Value = 3.65
$ws.Range("I7").Value = 4.7
$ws.Range("L7").Value = 1.26
$ws.Range("M7").Value = 3.15
$ws.Range("N7").Value = 1.78
$ws.Range("O7").Value = 1.83
$ws.Range("P7").Value = 1.39
$ws.Range("Q7").Value = 2.55
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 1.83
$ws.Range("T7").Value = 7
$ws.Range("U7").Value = 7.9
$ws.Range("V7").Value = 8
$ws.Range("W7").Value = 12.5
$ws.Range("X7").Value = 13
$ws.Range("Y7").Value = 25
$ws.Range("Z7").Value = 10.5
$ws.Range("AC8").Value = 90
$ws.Range("AF8").Value = 60
$ws.Range("AG8").Value = 25
$ws.Range("M8").Value = 3.5
$ws.Range("R8").Value = 1.93
$ws.Range("S8").Value = 1.7
$ws.Range("X8").Value = 11.25
$ws.Range("Z8").Value = 12
$ws.Range("AB9").Value = 17
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 4.05
$ws.Range("Z9").Value = 8.25
$ws.Range("AA10").Value = 6.9
$ws.Range("AB10").Value = 14
$ws.Range("AC10").Value = 60
$ws.Range("AD10").Value = 450
$ws.Range("AE10").Value = 8.25
$ws.Range("AF10").Value = 10.5
$ws.Range("AG10").Value = 8.75
$ws.Range("AH10").Value = 19.5
$ws.Range("AI10").Value = 16.5
$ws.Range("AJ10").Value = 26
$ws.Range("H10").Value = 3.55
$ws.Range("I10").Value = 2.12
$ws.Range("L10").Value = 1.25
$ws.Range("M10").Value = 3.2
$ws.Range("N10").Value = 1.75
$ws.Range("O10").Value = 1.85
$ws.Range("R10").Value = 1.65
$ws.Range("S10").Value = 1.98
$ws.Range("T10").Value = 10.25
$ws.Range("U10").Value = 16
$ws.Range("V10").Value = 10.75
$ws.Range("W10").Value = 35
$ws.Range("X10").Value = 25
$ws.Range("Y10").Value = 32
$ws.Range("Z10").Value = 11.5
$ws.Range("AB13").Value = 19.5
$ws.Range("AC13").Value = 150
$ws.Range("AE13").Value = 8
$ws.Range("AG13").Value = 14
$ws.Range("AH13").Value = 65
$ws.Range("AI13").Value = 45
$ws.Range("AJ13").Value = 65
$ws.Range("G13").Value = 2.1
$ws.Range("H13").Value = 2.87
$ws.Range("I13").Value = 3.9
$ws.Range("J13").Value = 1.13
$ws.Range("K13").Value = 5.1
$ws.Range("L13").Value = 1.57
$ws.Range("M13").Value = 2.27
$ws.Range("N13").Value = 2.65
$ws.Range("O13").Value = 1.42
$ws.Range("P13").Value = 1.57
$ws.Range("Q13").Value = 2.25
$ws.Range("R13").Value = 2.2
$ws.Range("S13").Value = 1.6
$ws.Range("T13").Value = 5.2
$ws.Range("V13").Value = 9.5
$ws.Range("W13").Value = 19.5
$ws.Range("X13").Value = 22
$ws.Range("Y13").Value = 45
$ws.Range("Z13").Value = 5.1
$ws.Range("AB19").Value = 15
$ws.Range("AC19").Value = 60
$ws.Range("AI19").Value = 11.75
$ws.Range("AJ19").Value = 22
$ws.Range("G19").Value = 4.9
$ws.Range("H19").Value = 4.1
$ws.Range("O19").Value = 2.15
$ws.Range("T19").Value = 16
$ws.Range("U19").Value = 30
$ws.Range("V19").Value = 15.5
$ws.Range("Y19").Value = 40
$ws.Range("AB25").Value = 15
$ws.Range("AD25").Value = 301
$ws.Range("AE25").Value = 9
$ws.Range("AG25").Value = 12
$ws.Range("H25").Value = 3.1
$ws.Range("K25").Value = 8.5
$ws.Range("N25").Value = 2.1
$ws.Range("O25").Value = 1.7
$ws.Range("P25").Value = 1.44
$ws.Range("Q25").Value = 2.63
$ws.Range("T25").Value = 7.5
$ws.Range("Y25").Value = 34
$ws.Range("Z25").Value = 8.5
$ws.Range("AA26").Value = 7.5
$ws.Range("AB26").Value = 18
$ws.Range("AC26").Value = 80
$ws.Range("AE26").Value = 14
$ws.Range("AF26").Value = 35
$ws.Range("AG26").Value = 19
$ws.Range("AH26").Value = 120
$ws.Range("AI26").Value = 65
$ws.Range("AJ26").Value = 65
$ws.Range("G26").Value = 1.35
$ws.Range("H26").Value = 4.3
$ws.Range("I26").Value = 7.1
$ws.Range("N26").Value = 1.72
$ws.Range("O26").Value = 1.88
$ws.Range("T26").Value = 5.5
$ws.Range("U26").Value = 5.2
$ws.Range("V26").Value = 7.2
$ws.Range("W26").Value = 7.1
$ws.Range("X26").Value = 9.75
$ws.Range("Y26").Value = 24
$ws.Range("Z26").Value = 10.75
$ws.Range("G27").Value = 1.83
$ws.Range("H27").Value = 3.8
$ws.Range("I27").Value = 3.9
$ws.Range("W27").Value = 15
$ws.Range("AA31").Value = 7.9
$ws.Range("AB31").Value = 14.5
$ws.Range("AC31").Value = 55
$ws.Range("AD31").Value = 350
$ws.Range("AE31").Value = 16
$ws.Range("AF31").Value = 30
$ws.Range("AI31").Value = 40
$ws.Range("AJ31").Value = 40
$ws.Range("G31").Value = 1.62
$ws.Range("H31").Value = 3.95
$ws.Range("I31").Value = 4.6
$ws.Range("L31").Value = 1.19
$ws.Range("M31").Value = 3.65
$ws.Range("N31").Value = 1.6
$ws.Range("O31").Value = 2.07
$ws.Range("R31").Value = 1.62
$ws.Range("S31").Value = 2.02
$ws.Range("T31").Value = 8.25
$ws.Range("U31").Value = 8.5
$ws.Range("X31").Value = 12
$ws.Range("Y31").Value = 22
$ws.Range("Z31").Value = 13.5
$ws.Range("AB32").Value = 14.5
$ws.Range("AC32").Value = 55
$ws.Range("AD32").Value = 350
$ws.Range("AE32").Value = 14.5
$ws.Range("AF32").Value = 26
$ws.Range("AG32").Value = 14
$ws.Range("AH32").Value = 65
$ws.Range("AI32").Value = 37
$ws.Range("AJ32").Value = 37
$ws.Range("G32").Value = 1.65
$ws.Range("I32").Value = 4.35
$ws.Range("K32").Value = 8.75
$ws.Range("M32").Value = 4
$ws.Range("N32").Value = 1.62
$ws.Range("O32").Value = 2.15
$ws.Range("P32").Value = 1.32
$ws.Range("Q32").Value = 3.1
$ws.Range("R32").Value = 1.65
$ws.Range("S32").Value = 2.1
$ws.Range("T32").Value = 8.5
$ws.Range("U32").Value = 8.75
$ws.Range("W32").Value = 13
$ws.Range("X32").Value = 12.5
$ws.Range("Y32").Value = 22
$ws.Range("Z32").Value = 8.75
$ws.Range("K33").Value = 12
$ws.Range("AG34").Value = 12
$ws.Range("H34").Value = 3.3
$ws.Range("J34").Value = 1.05
$ws.Range("K34").Value = 11
$ws.Range("N34").Value = 1.85
$ws.Range("O34").Value = 1.95
$ws.Range("V34").Value = 9.5
$ws.Range("AA36").Value = 6.1
$ws.Range("AB36").Value = 12.5
$ws.Range("AC36").Value = 55
$ws.Range("AD36").Value = 400
$ws.Range("AE36").Value = 9.5
$ws.Range("AF36").Value = 15
$ws.Range("AG36").Value = 10
$ws.Range("AH36").Value = 35
$ws.Range("AI36").Value = 23
$ws.Range("AJ36").Value = 29
$ws.Range("G36").Value = 2.42
$ws.Range("I36").Value = 2.8
$ws.Range("L36").Value = 1.28
$ws.Range("M36").Value = 3.05
$ws.Range("N36").Value = 1.83
$ws.Range("O36").Value = 1.78
$ws.Range("P36").Value = 1.4
$ws.Range("Q36").Value = 2.52
$ws.Range("R36").Value = 1.62
$ws.Range("S36").Value = 2.02
$ws.Range("T36").Value = 8.25
$ws.Range("U36").Value = 12.5
$ws.Range("V36").Value = 9.25
$ws.Range("W36").Value = 26
$ws.Range("X36").Value = 19.5
$ws.Range("Y36").Value = 27
$ws.Range("Z36").Value = 9.75
